$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Totals area (rows 1 and 4) ---
# Row 1: Total Hours formula now excludes row 301 (off-sheet) -> B300
$ws.Range("F1").Formula = "=SUM(B2,B3:B300)"

# Row 4: "Week 3 Hours" label + new weekly total formula (F4 didn't exist before)
$ws.Range("E4").Value = "Week 3 Hours"
$ws.Range("F4").Formula = "=SUM(B44:B60)"

# --- Week 3 log entries (rows 49-60) ---
# Rows 49-52 are untouched in content; only their "spans" bookkeeping shifts
# automatically once column D stops being used further down the sheet.

# Old row 53 ("Making the new GUI Map" placeholder with no hours/date) is
# effectively dropped, and every entry below it shifts up by one row. Row 53
# now carries what used to be row 54's entry, with its own hours/date filled in.
$ws.Range("C49").Copy($ws.Range("C53"))
$ws.Range("A53").Value = "Designing new Node Travel"
$ws.Range("B53").Value = 1
$ws.Range("C53").Value = 41931

# Row 54 now carries what used to be row 55's entry
$ws.Range("A54").Value = "Making new Prefabs, Switching to code graph based instead of GUI"
$ws.Range("B54").Value = 7

# Row 55 now carries what used to be row 56's entry
$ws.Range("A55").Value = "Fixing Errors in new setup"
$ws.Range("B55").Value = 2.5

# Row 56 now carries what used to be row 57's entry
$ws.Range("A56").Value = "Fixed teleporting lag/glitch"
$ws.Range("B56").Value = 2

# Row 57 now carries what used to be row 58's entry
$ws.Range("A57").Value = "Working on spawning Objects at end of the maze"
$ws.Range("B57").Value = 3
$ws.Range("C57").Value = 41932

# Row 58 now carries what used to be row 59's entry
$ws.Range("A58").Value = "Working on spawning Objects at end of the maze"
$ws.Range("B58").Value = 3.5

# Row 59 now carries what used to be row 60's entry (minus the stray note)
$ws.Range("A59").Value = "Working on spawning Objects at end of the maze"
$ws.Range("B59").Value = 3

# Row 60: new final entry for the week - bug fixing / cleanup
$ws.Range("A60").Value = "Fixing Bugs/Cleaning a bit of code"
$ws.Range("B60").Value = 1
$ws.Range("C60").Value = 41932
$ws.Range("D60").ClearContents()

# --- View state ---
$ws.Range("B21").Select()
